$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Password  changing" paragraph - remove the gramStart/gramEnd
# proofing-error markers and merge the two runs into a single run whose text
# is "Password  changing " (trailing space kept, xml:space="preserve").
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd() -eq "Password  changing") {
        $pwdRange = $para.Range
        $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="710BA3AC" w14:textId="15AA853F" w:rsidR="00624EEE" w:rsidRDefault="00624EEE" w:rsidP="00624EEE"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Password  changing </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $pwdRange.InsertXML($xml)
        break
    }
}

# ---------------------------------------------------------------------------
# Change 2: highlight the "logger" paragraph (paragraph mark + run) yellow.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd() -eq "logger") {
        $para.Range.Font.HighlightColorIndex = 7
        break
    }
}

# ---------------------------------------------------------------------------
# Change 3: mark the run that holds the inline picture as "do not spell
# check" (w:noProof) - this is the last paragraph of the document.
# ---------------------------------------------------------------------------
$shape = $d.InlineShapes.Item(1)
$shape.Range.Font.NoProofing = $true
